$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "Davide Rosà"
$ws.Range("B30").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C30").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("D30").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E30").Value = "Filippo Benetti | I Magnifici"
$ws.Range("F30").Value = "Mattia Tezzele | U.SGUARNA"
